$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 6452461
$ws.Range("I137").Value = 811.48
$ws.Range("J137").Value = 33334334
$ws.Range("K137").Value = 2434.44
$ws.Range("L137").Value = 100003002
$ws.Range("M137").Value = 115.5599999999999
$ws.Range("N137").Value = -100008102

$ws.Range("H138").Value = 1404.6731
$ws.Range("I138").Value = 1171.1277
$ws.Range("J138").Value = 3600
$ws.Range("K138").Value = 3513.3831
$ws.Range("L138").Value = 10800
$ws.Range("M138").Value = 1626.6169
$ws.Range("N138").Value = -21080

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7577013
$ws.Range("I61").Value = 9616502
$ws.Range("J61").Value = 1766.2142
$ws.Range("K61").Value = 9616502
$ws.Range("L61").Value = 1766.2142
$ws.Range("M61").Value = -9616290
$ws.Range("N61").Value = -2190.2142

$ws.Range("H74").Value = 13516269
$ws.Range("I74").Value = 23811412
$ws.Range("J74").Value = 3894.5
$ws.Range("K74").Value = 23811412
$ws.Range("L74").Value = 3894.5
$ws.Range("M74").Value = -23810538
$ws.Range("N74").Value = -5642.5

$ws.Range("H77").Value = 13516269
$ws.Range("I77").Value = 23811412
$ws.Range("J77").Value = 3894.5
$ws.Range("K77").Value = 119057060
$ws.Range("L77").Value = 19472.5
$ws.Range("M77").Value = -119052692
$ws.Range("N77").Value = -28208.5

$ws.Range("H122").Value = 15311.5
$ws.Range("I122").Value = 27873
$ws.Range("J122").Value = 2750
$ws.Range("K122").Value = 83619
$ws.Range("L122").Value = 8250
$ws.Range("M122").Value = -81169
$ws.Range("N122").Value = -13150

$ws.Range("H132").Value = 22731930
$ws.Range("I132").Value = 27781912
$ws.Range("J132").Value = 7007
$ws.Range("K132").Value = 83345736
$ws.Range("L132").Value = 21021
$ws.Range("M132").Value = -83343206
$ws.Range("N132").Value = -26081

$ws.Range("H136").Value = 7577013
$ws.Range("I136").Value = 9616502
$ws.Range("J136").Value = 1766.2142
$ws.Range("K136").Value = 28849506
$ws.Range("L136").Value = 5298.642599999999
$ws.Range("M136").Value = -28846956
$ws.Range("N136").Value = -10398.6426

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2804.1538
$ws.Range("I134").Value = 1989.6296
$ws.Range("J134").Value = 4636.8335
$ws.Range("K134").Value = 5968.8888
$ws.Range("L134").Value = 13910.5005
$ws.Range("M134").Value = -3433.8888
$ws.Range("N134").Value = -18980.5005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6293150
$ws.Range("I31").Value = 3834.7954
$ws.Range("J31").Value = 37040910
$ws.Range("K31").Value = 3834.7954
$ws.Range("L31").Value = 37040910
$ws.Range("M31").Value = -3539.7954
$ws.Range("N31").Value = -37041500

$ws.Range("H34").Value = 6293150
$ws.Range("I34").Value = 3834.7954
$ws.Range("J34").Value = 37040910
$ws.Range("K34").Value = 3834.7954
$ws.Range("L34").Value = 37040910
$ws.Range("M34").Value = -3632.7954
$ws.Range("N34").Value = -37041314

$ws.Range("H58").Value = 1767.3478
$ws.Range("I58").Value = 613.1053000000001
$ws.Range("J58").Value = 7250
$ws.Range("K58").Value = 613.1053000000001
$ws.Range("L58").Value = 7250
$ws.Range("M58").Value = -410.1053000000001
$ws.Range("N58").Value = -7656

$ws.Range("H136").Value = 1767.3478
$ws.Range("I136").Value = 613.1053000000001
$ws.Range("J136").Value = 7250
$ws.Range("K136").Value = 1839.3159
$ws.Range("L136").Value = 21750
$ws.Range("M136").Value = 710.6840999999999
$ws.Range("N136").Value = -26850

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 851
$ws.Range("I68").Value = 663.6667
$ws.Range("J68").Value = 907.2
$ws.Range("K68").Value = 1991.0001
$ws.Range("L68").Value = 2721.6
$ws.Range("M68").Value = -1180.0001
$ws.Range("N68").Value = -4343.6

$ws.Range("H71").Value = 851
$ws.Range("I71").Value = 663.6667
$ws.Range("J71").Value = 907.2
$ws.Range("K71").Value = 5973.0003
$ws.Range("L71").Value = 8164.8
$ws.Range("M71").Value = -1917.0003
$ws.Range("N71").Value = -16276.8

$ws.Range("H129").Value = 2227.1482
$ws.Range("I129").Value = 1048.6364
$ws.Range("J129").Value = 3037.375
$ws.Range("K129").Value = 3145.9092
$ws.Range("L129").Value = 9112.125
$ws.Range("M129").Value = 1854.0908
$ws.Range("N129").Value = -19112.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("N110").Value = 0
$ws.Range("L110").ClearContents()

$ws.Range("H113").Value = 101617.2
$ws.Range("I113").Value = 101617.2
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 101617.2
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = -99447.2
$ws.Range("M113").ClearContents()

$ws.Range("H114").Value = 60000
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 60000
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 60000
$ws.Range("N114").Value = -68678

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("N116").Value = 0
$ws.Range("L116").ClearContents()

$ws.Range("H117").Value = 59310
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 59310
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 59310
$ws.Range("N117").Value = -66194

$ws.Range("H122").Value = 11114944
$ws.Range("I122").Value = 16669416
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 50008248
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -50005798
$ws.Range("N122").Value = -22900

$ws.Range("H132").Value = 4097.1953
$ws.Range("I132").Value = 3164.5862
$ws.Range("J132").Value = 6351
$ws.Range("K132").Value = 9493.758600000001
$ws.Range("L132").Value = 19053
$ws.Range("M132").Value = -6963.758600000001
$ws.Range("N132").Value = -24113

$ws.Range("H140").Value = 69000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 69000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 69000
$ws.Range("N140").Value = -79360

$ws.Range("H141").Value = 400214.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 400214.5
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 400214.5
$ws.Range("N141").Value = -410574.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5656
$ws.Range("I7").Value = 6706.8335
$ws.Range("J7").Value = 5025.5
$ws.Range("K7").Value = 6706.8335
$ws.Range("L7").Value = 5025.5
$ws.Range("M7").Value = -6594.8335
$ws.Range("N7").Value = -5249.5

$ws.Range("H61").Value = 1637.4166
$ws.Range("I61").Value = 1522.1111
$ws.Range("J61").Value = 1983.3334
$ws.Range("K61").Value = 1522.1111
$ws.Range("L61").Value = 1983.3334
$ws.Range("M61").Value = -1320.1111
$ws.Range("N61").Value = -2387.3334

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("N110").Value = 0
$ws.Range("L110").ClearContents()

$ws.Range("H113").Value = 1637.4166
$ws.Range("I113").Value = 1522.1111
$ws.Range("J113").Value = 1983.3334
$ws.Range("K113").Value = 1522.1111
$ws.Range("L113").Value = 1983.3334
$ws.Range("M113").Value = 647.8888999999999
$ws.Range("N113").Value = -6323.3334

$ws.Range("H126").Value = 5656
$ws.Range("I126").Value = 6706.8335
$ws.Range("J126").Value = 5025.5
$ws.Range("K126").Value = 20120.5005
$ws.Range("L126").Value = 15076.5
$ws.Range("M126").Value = -17650.5005
$ws.Range("N126").Value = -20016.5

$ws.Range("H132").Value = 8070472.5
$ws.Range("I132").Value = 3125.4546
$ws.Range("J132").Value = 27790654
$ws.Range("K132").Value = 9376.363799999999
$ws.Range("L132").Value = 83371962
$ws.Range("M132").Value = -6846.363799999999
$ws.Range("N132").Value = -83377022

$ws.Range("H136").Value = 9620556
$ws.Range("I136").Value = 11365196
$ws.Range("J136").Value = 25038.125
$ws.Range("K136").Value = 34095588
$ws.Range("L136").Value = 75114.375
$ws.Range("M136").Value = -34093038
$ws.Range("N136").Value = -80214.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1806.7949
$ws.Range("I132").Value = 1441.9429
$ws.Range("J132").Value = 4999.25
$ws.Range("K132").Value = 4325.8287
$ws.Range("L132").Value = 14997.75
$ws.Range("M132").Value = -1795.8287
$ws.Range("N132").Value = -20057.75

$ws.Range("H136").Value = 1331.6552
$ws.Range("I136").Value = 964.72
$ws.Range("J136").Value = 3625
$ws.Range("K136").Value = 2894.16
$ws.Range("L136").Value = 10875
$ws.Range("M136").Value = -344.1599999999999
$ws.Range("N136").Value = -15975

$ws.Range("H140").Value = 113214.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 113214.5
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 113214.5
$ws.Range("N140").Value = -123574.5

$ws.Range("H141").Value = 65000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 65000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 65000
$ws.Range("N141").Value = -79360
